# NIT-9012644226.xlsx — "Elimina EC anteriores y se agregan nuevos,
# se modifica base de datos"
#
# The "Periodo Mora" column (E16:E21) held the six billing periods in
# ascending order (2501..2506). The new account statement replaces those
# with the same six periods listed in descending order (2506..2501), i.e.
# the most recent period first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2504"
$ws.Range("E19").Value = "2503"
$ws.Range("E20").Value = "2502"
$ws.Range("E21").Value = "2501"
